$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.232.06'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.651.02'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.17'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3902'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3890'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.001'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.376'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '49.25'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -6.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08513'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.52'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -6.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.197'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001295'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.558'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.651.38'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.13'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.36'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06932'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.993'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.90'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.246.41'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.373'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.772'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -6.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.64'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '158.60'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.699'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '143.47'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.348'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -13.77%  '
$ws.Range("E32").Value = '  -8.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.831.99'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.046'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08138'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02945'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2729'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09336'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.483'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.13'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7713'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.25'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.19'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.521'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6944'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.109'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.33%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08460'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.281'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -9.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.05'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.96%  '
